$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 150, pushing existing rows 150-223 down to 151-224.
$ws.Rows.Item(150).Insert()

# Populate the newly inserted row 150 with the new data record.
$ws.Cells.Item(150, 1).Value = 7
$ws.Cells.Item(150, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(150, 3).Value = "Ñuble"
$ws.Cells.Item(150, 4).Value = 45027
$ws.Cells.Item(150, 5).Value = 16
$ws.Cells.Item(150, 6).Value = "Fruta"
$ws.Cells.Item(150, 7).Value = 100101
$ws.Cells.Item(150, 8).Value = "Berries"
$ws.Cells.Item(150, 9).Value = 100101007
$ws.Cells.Item(150, 10).Value = "Kiwi"
$ws.Cells.Item(150, 11).Value = "Hayward"
$ws.Cells.Item(150, 12).Value = "Primera"
$ws.Cells.Item(150, 13).Value = 70
$ws.Cells.Item(150, 14).Value = 16000
$ws.Cells.Item(150, 15).Value = 16000
$ws.Cells.Item(150, 16).Value = 16000
$ws.Cells.Item(150, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(150, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(150, 19).Value = 889
$ws.Cells.Item(150, 20).Value = 18
